# gui added new panel, checkboxes panel and future goals
# Adds two new arrival rows (row 21 and row 22) to the "Main Data" sheet,
# continuing the Friday, Jan 13 entries already present in rows 19-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: FR5217 from Dublin (DUB), Ryanair B738 (EI-EXE)
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(21, 3).Value = "6:15 PM"
$ws.Cells.Item(21, 4).Value = "FR5217"
$ws.Cells.Item(21, 5).Value = "Dublin"
$ws.Cells.Item(21, 6).Value = "(DUB)"
$ws.Cells.Item(21, 7).Value = "Ryanair "
$ws.Cells.Item(21, 8).Value = "B738"
$ws.Cells.Item(21, 9).Value = "(EI-EXE)"
$ws.Cells.Item(21, 10).Value = "5:57 PM"
$ws.Cells.Item(21, 11).ClearFormats()
$ws.Cells.Item(21, 12).Value = "0 hours, -18 minutes"
$ws.Cells.Item(21, 13).ClearFormats()

# Row 22: FR6111 from Gdansk (GDN), Ryanair B738 (SP-RSW)
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(22, 3).Value = "8:55 PM"
$ws.Cells.Item(22, 4).Value = "FR6111"
$ws.Cells.Item(22, 5).Value = "Gdansk"
$ws.Cells.Item(22, 6).Value = "(GDN)"
$ws.Cells.Item(22, 7).Value = "Ryanair "
$ws.Cells.Item(22, 8).Value = "B738"
$ws.Cells.Item(22, 9).Value = "(SP-RSW)"
$ws.Cells.Item(22, 10).Value = "8:57 PM"
$ws.Cells.Item(22, 11).ClearFormats()
$ws.Cells.Item(22, 12).Value = "0 hours, 2 minutes"
$ws.Cells.Item(22, 13).ClearFormats()
